$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'63.810.58"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -3.83%  '
$ws.Range('D3').Value = "'3.102.16"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -5.07%  '
$ws.Range('D5').Value = "'608.27"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.08%  '
$ws.Range('D6').Value = "'144.89"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -8.17%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').Value = "'3.099.06"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -5.14%  '
$ws.Range('D9').Value = "'0.519"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -4.69%  '
$ws.Range('D10').Value = "'0.150"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -7.26%  '
$ws.Range('D11').Value = "'5.23"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -9.68%  '
$ws.Range('D12').Value = "'0.468"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -5.69%  '
$ws.Range('D13').Value = "'0.0000249"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -7.87%  '
$ws.Range('D14').Value = "'35.19"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -9.98%  '
$ws.Range('D15').Value = "'3.608.08"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -5.19%  '
$ws.Range('E16').Value = '  +0.85%  '
$ws.Range('D17').Value = "'63.853.79"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -3.87%  '
$ws.Range('D18').Value = "'3.095.53"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -5.33%  '
$ws.Range('D19').Value = "'6.81"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -8.17%  '
$ws.Range('D20').Value = "'475.55"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -5.80%  '
$ws.Range('D21').Value = "'14.62"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -5.85%  '
$ws.Range('D22').Value = "'0.698"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -7.56%  '
$ws.Range('D23').Value = "'7.69"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -5.58%  '
$ws.Range('D24').Value = "'13.53"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -7.77%  '
$ws.Range('D25').Value = "'83.32"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -4.27%  '
$ws.Range('D26').Value = "'1.00"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.03%  '
$ws.Range('D27').Value = "'2.77"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -9.19%  '
$ws.Range('D28').Value = "'8.38"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -9.35%  '
$ws.Range('E29').Value = '  -10.94%  '
$ws.Range('D30').Value = "'6.68"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -5.37%  '
$ws.Range('D31').Value = "'0.113"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -12.49%  '
$ws.Range('E32').Value = '  +0.03%  '
$ws.Range('D33').Value = "'2.71"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -5.75%  '
$ws.Range('D34').Value = "'26.14"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -6.69%  '
$ws.Range('E35').Value = '  -3.96%  '
$ws.Range('D36').Value = "'5.92"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -8.70%  '
$ws.Range('D37').Value = "'52.64"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -5.37%  '
$ws.Range('D38').Value = "'0.0₃0737"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -6.91%  '
$ws.Range('D39').Value = "'459.52"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -7.10%  '
$ws.Range('D40').Value = "'2.94"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -14.05%  '
$ws.Range('D41').Value = "'0.0391"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -7.96%  '
$ws.Range('E42').Value = '  -8.19%  '
$ws.Range('D43').Value = "'8.33"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -5.65%  '
$ws.Range('D44').Value = "'2.831.13"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -6.26%  '
$ws.Range('D45').Value = "'0.266"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -9.22%  '
$ws.Range('D46').Value = "'2.25"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -11.58%  '
$ws.Range('E47').Value = '  -3.74%  '
$ws.Range('E48').Value = '  -0.01%  '
$ws.Range('D49').Value = "'26.11"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -9.76%  '
$ws.Range('D50').Value = "'0.113"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -5.45%  '
$ws.Range('D51').Value = "'118.34"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.58%  '
